$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 16.69780766666667
$ws.Range("N2").Value = 50.093423
$ws.Range("O2").Value = 0.3087131902856588
$ws.Range("P2").Value = 0.3087131902856588
$ws.Range("Q2").Value = 0.777600205229
$ws.Range("R2").Value = 6.998401847061
$ws.Range("S2").Value = 0.3087131902856588
$ws.Range("T2").Value = 0.3087131902856588

# Row 3
$ws.Range("N3").Value = 83.628069
$ws.Range("O3").Value = 0.5153787949052554
$ws.Range("P3").Value = 0.5153787949052554
$ws.Range("S3").Value = 0.5153787949052554
$ws.Range("T3").Value = 0.5153787949052554

# Row 4
$ws.Range("O4").Value = 0.1759080148090858
$ws.Range("P4").Value = 0.1759080148090857
$ws.Range("S4").Value = 0.1759080148090858
$ws.Range("T4").Value = 0.1759080148090857
